$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.913.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4670'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3690'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07363'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8717'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.40'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.815.20'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.380'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07072'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.515'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008700'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.931.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.322'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.022.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.70%  '
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.46'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.173'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.341'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08922'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7695'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.506'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.904'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  -3.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01963'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05289'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.933'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.261'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5327'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.354'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.97%  '
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.443'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4932'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.673'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.000'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.86'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06297'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
